$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 16 (pushes totals row 16->17, footer row 17->18)
$ws.Rows("16:16").Insert()

# Populate the new row 16 by copying row 15's values + formatting (reuses the same
# cell styles as the rows above instead of inventing new ones), then we will
# rewrite rows 13-15 below so the product list ends up in the right order.
$ws.Range("A15:N15").Copy($ws.Range("A16:N16"))
$ws.Range("A16").Value = 13

# Row 15 becomes what used to be row 14's product (ZURCAL)
$ws.Range("B15").Value = "ZURCAL 40MG 14 GASTRO RESISTANT TAB"
$ws.Range("H15").Value = "4:0"
$ws.Range("L15").Value = 96
$ws.Range("N15").Value = "1:0"

# Row 14 becomes what used to be row 13's product (VOLTAREN)
$ws.Range("B14").Value = "VOLTAREN 75MG/3ML 3 AMP."
$ws.Range("H14").Value = "5:1"
$ws.Range("L14").Value = 17
$ws.Range("N14").Value = "0:0"

# Row 13 becomes the newly added product (TRICOVEL)
$ws.Range("B13").Value = "TRICOVEL 30 TABS."
$ws.Range("H13").Value = "0:0"
$ws.Range("L13").Value = 531
$ws.Range("N13").Value = "1:0"

# Update the grand total (was 690, now 690 + 531 = 1221)
$ws.Range("K17").Value = 1221
